# "gathercolums" style edit: widen the grouped-columns block by 4 columns.
#
# Originally columns C..BF held, for every person row, the same repeated
# "group" marker (e.g. "group1"/"group2"/"group3"), with the very last
# column (BF) doing double duty: sometimes it simply continued the
# repeated marker, sometimes it instead carried a trailing numeric tally
# (rows 2-5) or was left blank (rows 6-15).
#
# After the change the repeated block is 4 columns wider (it now goes all
# the way to BI) and the old trailing value that used to live in BF is
# pushed out to the new last column, BJ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 15
$groupCol = 3          # column C: holds the per-row group marker
$oldLastCol = 58        # column BF, before the edit
$newCols = @(58, 59, 60, 61)  # BF, BG, BH, BI after the edit
$newLastCol = 62        # column BJ, after the edit

# Step 1: insert 4 new columns where BF used to be. Excel/the engine
# shifts the previous BF column (and anything right of it) four places
# to the right, so the old BF content ends up in BJ, and the newly
# inserted BF:BI cells inherit the same cell style as the column being
# pushed aside.
$insertRange = $ws.Range("BF1:BI" + $lastDataRow)
$insertRange.Insert(-4161)

# Step 2: the insert above also shifts the width/style metadata that was
# defined, far outside the used data area, for columns 843-1024 -- that
# metadata is unrelated to this edit and must stay put. Deleting 4 blank
# columns well to the right of the real data (and to the left of that
# metadata block) shifts it back left by 4, exactly undoing that
# incidental side effect while leaving every real cell untouched.
$farRange = $ws.Range($ws.Cells.Item(1, 100), $ws.Cells.Item(1, 103)).EntireColumn
$farRange.Delete(-4161)

# Step 3: fill the newly available BF:BI columns with the same "group"
# marker that already fills columns C through BE on each row. Rows with
# nothing in column C (the blank row 10) simply receive blank cells,
# matching the existing pattern for that row.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $groupValue = $ws.Cells.Item($r, $groupCol).Value()
    $fillRange = $ws.Range($ws.Cells.Item($r, $newCols[0]), $ws.Cells.Item($r, $newCols[3]))
    $fillRange.Value = $groupValue
}
